$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-02 Wednesday" "2024-10-03 Thursday"

Replace-Text "214÷2=" "772÷7="
Replace-Text "161÷2=" "775÷2="
Replace-Text "705÷5=" "915÷4="
Replace-Text "456÷5=" "737÷9="
Replace-Text "384÷8=" "761÷4="
Replace-Text "544÷5=" "963÷2="
Replace-Text "112÷8=" "431÷8="
Replace-Text "627÷7=" "183÷9="
Replace-Text "391÷3=" "503÷8="
Replace-Text "878÷9=" "623÷5="
Replace-Text "509÷8=" "925÷8="
Replace-Text "973÷7=" "680÷5="
Replace-Text "250÷8=" "131÷9="
Replace-Text "465÷4=" "934÷5="
Replace-Text "174÷9=" "499÷5="
Replace-Text "425÷4=" "348÷6="
Replace-Text "721÷2=" "717÷3="
Replace-Text "539÷6=" "482÷3="
Replace-Text "623÷4=" "584÷5="
Replace-Text "409÷4=" "141÷2="
Replace-Text "669÷3=" "599÷2="
Replace-Text "118÷9=" "389÷3="
Replace-Text "470÷7=" "987÷9="
Replace-Text "868÷6=" "729÷2="
Replace-Text "887÷4=" "187÷8="
